$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 274, pushing the existing rows 274-339
# down to 275-340 (values/formatting travel with them automatically).
$ws.Rows(274).Insert()

# Populate the newly inserted row 274 with the new weekly data point.
$ws.Cells.Item(274, 1).Value = 3
$ws.Cells.Item(274, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(274, 3).Value = "Coquimbo"
$ws.Cells.Item(274, 4).Value = 44754
$ws.Cells.Item(274, 5).Value = 5
$ws.Cells.Item(274, 6).Value = 100112039
$ws.Cells.Item(274, 7).Value = "Ciboulette"
$ws.Cells.Item(274, 8).Value = "Sin especificar"
$ws.Cells.Item(274, 9).Value = "Primera"
$ws.Cells.Item(274, 10).Value = 110
$ws.Cells.Item(274, 11).Value = 1500
$ws.Cells.Item(274, 12).Value = 1500
$ws.Cells.Item(274, 13).Value = 1500
$ws.Cells.Item(274, 14).Value = "`$/docena de atados"
$ws.Cells.Item(274, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(274, 16).Value = 500
$ws.Cells.Item(274, 17).Value = 3
$ws.Cells.Item(274, 18).Value = "Hortaliza"
